$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.308.93"
$ws.Range("E2").Value = "  +1.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5264"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.009"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2688"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06462"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07516"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.39%  "

$ws.Range("D12").Value = "1.696.64"
$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.515"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.38%  "

$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008496"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.84%  "

$ws.Range("D17").Value = "26.320.95"
$ws.Range("E17").Value = "  +0.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.917"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.63"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.191"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.009"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.766"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.18%  "

$ws.Range("E26").Value = "  +5.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06476"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.362"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.18%  "

$ws.Range("E30").Value = "  +0.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.588"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.582"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.656"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.027"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6200"
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.405"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.742"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.288"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.77%  "

$ws.Range("D39").Value = "1.115.96"
$ws.Range("E39").Value = "  +3.84%  "

$ws.Range("E40").Value = "  +0.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8720"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.55%  "

$ws.Range("E42").Value = "  +0.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").Value = "1.829.93"
$ws.Range("E44").Value = "  +0.91%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.10%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.46%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.175"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.77%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05268"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4298"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.074"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.38%  "
